$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full corrected results table (70 data rows, header in row 1 untouched).
# The "weather" dataset block and the previous partial "IMDB reviews" block
# have been replaced by a freshly uploaded, complete "IMDB reviews" block,
# and values throughout (notably the "adult" dataset rows, which had used
# the wrong target column) were corrected.
$data = @(
    @(0.990159901599016, 0.998769987699877, 0.993849938499385, 0.996309963099631, "mushrooms"),
    @(0.990159901599016, 1, 0.985239852398524, 0.982779827798278, "mushrooms"),
    @(1, 1, 1, 0.988929889298893, "mushrooms"),
    @(0.995079950799508, 1, 0.9864698646986469, 0.9876998769987699, "mushrooms"),
    @(0.9950738916256158, 0.9987684729064039, 0.9926108374384236, 0.9876847290640394, "mushrooms"),
    @(0.9963054187192119, 1, 0.9950738916256158, 0.9815270935960592, "mushrooms"),
    @(0.9913793103448276, 0.9987684729064039, 0.9938423645320197, 0.9901477832512315, "mushrooms"),
    @(0.9667487684729064, 1, 0.9926108374384236, 0.9802955665024631, "mushrooms"),
    @(0.9938423645320197, 0.9963054187192119, 0.9938423645320197, 0.9827586206896551, "mushrooms"),
    @(0.9938423645320197, 0.9987684729064039, 0.9938423645320197, 0.9913793103448276, "mushrooms"),
    @(0.8580588105239885, 0.8589431793057705, 0.8633650232146806, 0.863143931019235, "adult"),
    @(0.8651337607782446, 0.8671235905372541, 0.8686712359053725, 0.8695556046871545, "adult"),
    @(0.8615656789031402, 0.8637770897832817, 0.8642193719593101, 0.8657673595754091, "adult"),
    @(0.8710747456877488, 0.8761609907120743, 0.8752764263600177, 0.8728438743918621, "adult"),
    @(0.8586908447589562, 0.8622291021671826, 0.865546218487395, 0.862671384343211, "adult"),
    @(0.8538257408226448, 0.8646616541353384, 0.862671384343211, 0.862671384343211, "adult"),
    @(0.8664307828394515, 0.8653250773993808, 0.8701901813356921, 0.8659885006634233, "adult"),
    @(0.865546218487395, 0.871517027863777, 0.8730650154798761, 0.8693056169836355, "adult"),
    @(0.8704113224237063, 0.8719593100398054, 0.8759398496240601, 0.8706324635117205, "adult"),
    @(0.8675364882795223, 0.8681999115435648, 0.8759398496240601, 0.8726227333038479, "adult"),
    @(0.8, 0.8085106382978723, 0.8028368794326242, 0.8070921985815603, "churn"),
    @(0.8042553191489362, 0.7929078014184398, 0.8141843971631205, 0.8042553191489362, "churn"),
    @(0.7971631205673759, 0.7843971631205674, 0.7815602836879433, 0.7787234042553192, "churn"),
    @(0.7826704545454546, 0.7826704545454546, 0.78125, 0.7897727272727273, "churn"),
    @(0.8068181818181818, 0.7741477272727273, 0.7911931818181818, 0.7940340909090909, "churn"),
    @(0.8068181818181818, 0.7755681818181818, 0.7897727272727273, 0.7755681818181818, "churn"),
    @(0.8025568181818182, 0.7883522727272727, 0.7926136363636364, 0.7911931818181818, "churn"),
    @(0.8082386363636364, 0.7613636363636364, 0.7954545454545454, 0.7940340909090909, "churn"),
    @(0.7954545454545454, 0.7769886363636364, 0.7755681818181818, 0.8039772727272727, "churn"),
    @(0.7954545454545454, 0.7727272727272727, 0.7741477272727273, 0.7883522727272727, "churn"),
    @(1, 1, 1, 1, "credit card"),
    @(1, 1, 0.998, 1, "credit card"),
    @(0.998, 0.999, 0.998, 0.998, "credit card"),
    @(0.999, 1, 1, 0.999, "credit card"),
    @(1, 1, 0.989, 1, "credit card"),
    @(1, 0.999, 0.992, 0.999, "credit card"),
    @(1, 1, 1, 1, "credit card"),
    @(1, 1, 0.998, 1, "credit card"),
    @(0.999, 0.999, 0.997, 0.999, "credit card"),
    @(0.999, 0.999, 0.998, 0.999, "credit card"),
    @(1, 1, 1, 1, "prostate"),
    @(1, 1, 1, 1, "prostate"),
    @(0.9, 1, 1, 1, "prostate"),
    @(0.8, 0.8, 0.8, 0.8, "prostate"),
    @(0.9, 0.9, 0.9, 0.9, "prostate"),
    @(1, 1, 1, 1, "prostate"),
    @(0.9, 0.9, 0.9, 0.9, "prostate"),
    @(0.9, 1, 1, 1, "prostate"),
    @(0.8, 0.9, 1, 0.9, "prostate"),
    @(0.9, 0.9, 0.9, 0.9, "prostate"),
    @(0.875, 0.875, 1, 1, "leukemia"),
    @(0.875, 0.875, 0.875, 0.875, "leukemia"),
    @(0.7142857142857143, 1, 1, 1, "leukemia"),
    @(0.7142857142857143, 0.7142857142857143, 0.8571428571428571, 0.8571428571428571, "leukemia"),
    @(0.7142857142857143, 1, 1, 1, "leukemia"),
    @(1, 1, 1, 1, "leukemia"),
    @(1, 1, 1, 1, "leukemia"),
    @(0.7142857142857143, 1, 1, 1, "leukemia"),
    @(1, 1, 1, 1, "leukemia"),
    @(1, 1, 1, 1, "leukemia"),
    @(0.8149999999999999, 0.79, 0.805, 0.85, "IMDB reviews"),
    @(0.8149999999999999, 0.805, 0.835, 0.8100000000000001, "IMDB reviews"),
    @(0.76, 0.77, 0.775, 0.76, "IMDB reviews"),
    @(0.8149999999999999, 0.8100000000000001, 0.8149999999999999, 0.82, "IMDB reviews"),
    @(0.8, 0.765, 0.805, 0.82, "IMDB reviews"),
    @(0.805, 0.825, 0.805, 0.825, "IMDB reviews"),
    @(0.795, 0.775, 0.795, 0.785, "IMDB reviews"),
    @(0.8149999999999999, 0.835, 0.785, 0.83, "IMDB reviews"),
    @(0.83, 0.84, 0.85, 0.875, "IMDB reviews"),
    @(0.8, 0.765, 0.79, 0.805, "IMDB reviews")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 2 + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Original sheet had 80 data rows (through row 81); the corrected data only
# needs 70 (through row 71), so drop the now-unused trailing rows.
$lastRow = 1 + $data.Count
$ws.Rows("72:81").Delete()

Write-Host "Wrote $($data.Count) data rows; sheet now ends at row $lastRow."
